# Fruta / hortaliza, semanal
# Insert a new weekly record as row 28 (shifting the existing rows 28-85 down to 29-86).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 28; everything below shifts down by one.
$ws.Rows("28:28").Insert()

# Populate the newly inserted row 28 with the new record's data.
$ws.Range("A28").Value = 5
$ws.Range("B28").Value = "Macroferia Regional de Talca"
$ws.Range("C28").Value = "Maule"
$ws.Range("D28").Value = 44883
$ws.Range("E28").Value = 7
$ws.Range("F28").Value = 300000000
$ws.Range("G28").Value = "Espárragos"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 3000
$ws.Range("K28").Value = 1000
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = 1000
$ws.Range("N28").Value = '$/atado'
$ws.Range("O28").Value = "Provincia de Linares"
$ws.Range("P28").Value = 1000
$ws.Range("Q28").Value = 1
$ws.Range("R28").Value = "Hortaliza"

# Make sure the date cell keeps the same date/time number format used throughout column D.
$ws.Range("D28").NumberFormat = "YYYY-MM-DD HH:MM:SS"
